$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.980.97"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5096"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2582"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06358"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.82"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07770"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.279"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.633.95"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5471"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7756"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.31"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.995.37"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "196.52"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.427"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.931"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.087"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.890"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.35"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1234"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +7.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.869"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.63"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.78%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04878"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.274"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.221"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.41%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.374"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9141"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.571"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5546"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.090.51"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01572"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.001"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.531"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.597"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8053"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.15"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₈121"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.775.66"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.56%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.14%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.45"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05220"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.539"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.68%  "

